# Apply data repull / push-all-data / mean-calculation corrections
# to the "dSF" column (column F) for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F7").Value  = 3
$ws.Range("F9").Value  = 1
$ws.Range("F15").Value = 8
$ws.Range("F18").Value = 1
$ws.Range("F23").Value = -2
$ws.Range("F24").Value = -2
